$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.311.63"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.680.56"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.06"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5416"
$ws.Range("E6").Value = "  +6.00%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2692"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06474"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.00"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07547"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.532"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "1.674.34"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5794"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.82"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "26.316.57"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.912"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.89"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.95"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.218"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "146.18"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1300"
$ws.Range("E25").Value = "  +8.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.837"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.77"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06475"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.394"
$ws.Range("E29").Value = "  +4.18%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.582"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.577"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.667"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6160"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.399"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.244"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "1.111.83"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8720"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.60"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "1.829.74"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.20"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.172"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4288"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.077"
$ws.Range("E51").Value = "  +0.64%  "